# Add a "selection_type" column to the selection-control sheet, between the
# existing "step_type" column (B) and "stock_concentration" column (old C).
# Also flips which sheet is the active/selected tab: selection-control
# becomes active (was selection-step_generation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("selection-control")
$ws2 = $wb.Worksheets.Item("selection-step_generation")

# Insert a new column before the old column C (stock_concentration), shifting
# stock_concentration..step_type (and their header-row cell comments) right
# by one column. Inserting the column itself only shifts cell VALUES - any
# existing cell comments stay anchored to their old column letter - so the
# comments need to be relocated by hand afterwards.
$ws.Columns.Item(3).Insert()

function Move-CellComment($sheet, $fromRef, $toRef) {
    $fromCell = $sheet.Range($fromRef)
    $comment = $fromCell.Comment
    if ($comment -ne $null) {
        $text = $comment.Text()
        $comment.Delete()
        $sheet.Range($toRef).AddComment($text)
    }
}

# Walk right-to-left so a cell's old comment is relocated before it becomes
# the destination of the next shift.
Move-CellComment $ws "L1" "M1"
Move-CellComment $ws "K1" "L1"
Move-CellComment $ws "J1" "K1"
Move-CellComment $ws "I1" "J1"
Move-CellComment $ws "H1" "I1"
Move-CellComment $ws "G1" "H1"
Move-CellComment $ws "F1" "G1"
Move-CellComment $ws "E1" "F1"
Move-CellComment $ws "C1" "D1"

# Header + comment for the new column.
$ws.Range("C1").Value = "selection_type"
$ws.Range("C1").AddComment("chemical: using a chemical selection. Slow pumps will be used to adjust chemical concentration`nnon-chemical: anything that does not require fluidics")

# Row 2 (vial 0, AUTO) is non-chemical; all the rest are chemical.
$ws.Range("C2").Value = "non-chemical"
for ($r = 3; $r -le 17; $r++) {
    $ws.Range("C$r").Value = "chemical"
}

# Update selection / active-tab state: selection-control becomes the active
# sheet (tabSelected), selection-step_generation is no longer the active tab.
# Select on the sheet that should stay inactive FIRST (it keeps its own
# selection memory but loses "active" status), then select on the sheet
# that should end up active/tabSelected LAST.
$ws2.Range("E3").Select()
$ws.Range("D10").Select()
